# The commit adds a new "5 steps of a SCM plan" section after the existing
# intro paragraph. Every new paragraph (and the intro paragraph's own mark)
# picks up explicit Arial / color 202124 / white-shading run formatting that
# isn't reachable through the simple Range.Font / InsertAfter surface (it has
# to land on <w:pPr><w:rPr> for empty paragraphs), so we drive the whole
# rewrite through Range.InsertXML with a WordProcessingML "Flat OPC" package,
# which is the standard COM/VBA technique for precise OOXML-level inserts.
#
# Range.InsertXML replaces the contents of the range it's called on, so
# calling it on $d.Content swaps out the whole body while Word keeps the
# trailing <w:sectPr> (and every other part) untouched.

$d = $word.ActiveDocument

$wordOpenXml = @'
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"><w:body>
<w:p w14:paraId="63C1C82F" w14:textId="74413653" w:rsidR="00994139" w:rsidRDefault="00857A27" w:rsidP="00857A27">
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Configuration Management (CM) is the ongoing process of identifying and managing changes to deliverables and other work products. The Configuration Management Plan (CM Plan) is&#160;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>developed to define, document, control, implement, account for, and audit changes to the various components of this project</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>The 5 steps of a SCM plan</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Planning and Identification. The first step in the process is planning and identification. ...</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Version Control and Baseline. ...</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Change Control. ...</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Configuration Status Accounting. ...</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>Audits and Reviews.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="202124"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>CM is useful for companies to manage the various versions of their codes.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$d.Content.InsertXML($wordOpenXml)
